$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 9 for the new "Approve/Reject" use-case pair,
# pushing the existing "View property details" / "Feed backs" rows (and
# everything after them) down by two.
$ws.Rows("9:10").Insert()
$ws.Range("C9").Value = "view buyers/ sellers request"
$ws.Range("G9").Value = "view the details of requested peoples"
$ws.Range("C10").Value = "Approve/Reject request"
$ws.Range("G10").Value = "Admin can delete or confirm the request"

# The old combined "Use(Buyer/Seller)" / "Registration" row (now row 14
# after the shift above) needs to be split across two rows: "Use(Buyer/
# Seller)" stays alone on row 14, "Registration" moves to a new row 15.
$ws.Rows("15:15").Insert()
$ws.Range("C15").Value = $ws.Range("C14").Value2
$ws.Range("C14").Value = ""

# Give the final "Visitor" row (now row 25) the taller custom row height
# seen in the target sheet, and update the active selection to match.
$ws.Rows("25:25").RowHeight = 26.25
$ws.Range("A25:O26").Select()
